$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsZhCn.Range("D4").Value = "2016-02-22 09:03:20"
$wsDeDe.Range("D4").Value = "2016-02-22 09:03:32"
